$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 295
$ws.Range("I2").Value = 703
$ws.Range("J2").Value = 3069
$ws.Range("L2").Value = 780
$ws.Range("M2").Value = 54
$ws.Range("N2").Value = 530
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 36
$ws.Range("S2").Value = 345
$ws.Range("T2").Value = 554
$ws.Range("U2").Value = 48
$ws.Range("V2").Value = 4565
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 4694
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 66
$ws.Range("AA2").Value = 26
